$d = $word.ActiveDocument

# --- Locate the two paragraphs we need to touch -----------------------
# Paragraph 6: the "Data yang Digunakan" heading-style paragraph that
#              currently also hosts the _GoBack bookmark.
# Paragraph 7: the paragraph that starts the "Data yang digunakan..."
#              sentence and currently ends with "), di mana ".
$pHeading = $d.Paragraphs.Item(6)
$pBody = $d.Paragraphs.Item(7)

# --- 1) Extend paragraph 7 with the rest of the sentence --------------
$r = $pBody.Range.Duplicate
$r.Collapse(0)
$r.Font.Italic = 0
$r.Font.ItalicBi = 0
$r.InsertAfter("nilai dari setiap parameter dipisahkan oleh sebuah karakter ")

$r.Collapse(0)
$r.InsertAfter("semicolon")
$r.Font.Italic = 1
$r.Font.ItalicBi = 1

$r.Collapse(0)
$r.Font.Italic = 0
$r.Font.ItalicBi = 0
$r.InsertAfter(" (titik koma), seperti yang ditunjukkan pada gambar 3.1. Data ini didapat dari penelitian yang dilakukan oleh Rahmat ")

$r.Collapse(0)
$r.InsertAfter("et al.")
$r.Font.Italic = 1
$r.Font.ItalicBi = 1

$r.Collapse(0)
$r.Font.Italic = 0
$r.Font.ItalicBi = 0
$r.InsertAfter(" (2016)")

# --- 2) Insert a new (initially empty) paragraph after paragraph 7 ----
$pBody2 = $d.Paragraphs.Item(7)
$rEnd = $pBody2.Range.Duplicate
$rEnd.Collapse(0)
$rEnd.InsertParagraphAfter()

# --- 3) Move the _GoBack bookmark from the heading paragraph to the ---
#        brand-new empty paragraph. A bookmark collapsed on the sole
#        character (the paragraph mark) of an empty paragraph can't be
#        targeted directly, so stash a temporary placeholder character,
#        anchor the bookmark next to it, then remove the placeholder
#        (the bookmark sticks to its position).
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

$pNew = $d.Paragraphs.Item(8)
$rPlaceholder = $pNew.Range.Duplicate
$rPlaceholder.Collapse(1)
$rPlaceholder.InsertBefore("X")

$pNew2 = $d.Paragraphs.Item(8)
$rBookmark = $d.Range($pNew2.Range.Start, $pNew2.Range.Start)
$d.Bookmarks.Add("_GoBack", $rBookmark)

$rRemove = $d.Range($rBookmark.Start, $rBookmark.Start + 1)
$rRemove.Delete()
